# Sample_Data.xlsx: add a 'New Id' demo row behaviour by updating the
# example Operation values and the Mult Qty test description/min value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("Testing Advanced Quantities"): Operation changes from Change -> Add
$ws.Range("A2").Value = "Add"

# Row 5 ("Testing Multiple Quantity"): Operation changes from Skip -> Change,
# the brief description now mentions starting at 100, and Min Qty becomes 100.
$ws.Range("A5").Value = "Change"
$ws.Range("F5").Value = "Order in multiples of 5, starting with 100, to a maximum of 1000."
$ws.Range("Y5").Value = 100

# Leave the selection on G2, matching where the cursor ended up after the edit.
$ws.Range("G2").Select()
